# Conversation Quality Check page: sample Q&A rows added to the template.
# - A3 is updated to a new prompt ("你是哪个模型厂商提供的")
# - A4 / A5 are new prompt rows
# - Selection moves to D10 (matches the author's last cursor position)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "你是哪个模型厂商提供的"
$ws.Range("A4").Value = "只切一刀，如何把四个橘子分给四个小朋友？"
$ws.Range("A5").Value = "我拿水兑水，得到的是稀水还是浓水？"

[void]$ws.Range("D10").Select()
